$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.653.23'
$ws.Range('E2').Value = '  +1.25%  '

$ws.Range('D3').Value = '1.827.78'
$ws.Range('E3').Value = '  +2.00%  '

$ws.Range('E4').Value = '  +0.45%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.18'
$ws.Range('E5').Value = '  +0.85%  '

$ws.Range('E6').Value = '  +0.40%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4665'
$ws.Range('E7').Value = '  +3.48%  '

$ws.Range('E8').Value = '  +0.11%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07139'
$ws.Range('E9').Value = '  +0.91%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9040'
$ws.Range('E10').Value = '  +2.27%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07689'
$ws.Range('E11').Value = '  -0.59%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '19.45'
$ws.Range('E12').Value = '  -0.02%  '

$ws.Range('D13').Value = '1.853.57'
$ws.Range('E13').Value = '  +4.04%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.269'
$ws.Range('E14').Value = '  -0.17%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.378'
$ws.Range('E15').Value = '  +0.97%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '87.72'
$ws.Range('E16').Value = '  +3.36%  '

$ws.Range('E17').Value = '  +0.42%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008563'
$ws.Range('E18').Value = '  +0.69%  '

$ws.Range('E19').Value = '  +0.36%  '

$ws.Range('D20').Value = '26.693.86'
$ws.Range('E20').Value = '  +1.27%  '

$ws.Range('E21').Value = '  -0.18%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.027'
$ws.Range('E22').Value = '  +1.23%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.55'
$ws.Range('E23').Value = '  +0.08%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.909'

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.95'
$ws.Range('E25').Value = '  +1.28%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '17.95'
$ws.Range('E26').Value = '  +0.83%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.996'
$ws.Range('E27').Value = '  -1.67%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '113.78'
$ws.Range('E28').Value = '  +1.93%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.873'
$ws.Range('E29').Value = '  +0.74%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08827'
$ws.Range('E30').Value = '  +1.72%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.144'
$ws.Range('E31').Value = '  +2.40%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.848'
$ws.Range('E32').Value = '  +3.07%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.166'
$ws.Range('E33').Value = '  +5.88%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7368'
$ws.Range('E34').Value = '  +2.03%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.436'
$ws.Range('E35').Value = '  -0.06%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.082'
$ws.Range('E36').Value = '  +1.47%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01933'
$ws.Range('E37').Value = '  +0.41%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05163'
$ws.Range('E38').Value = '  +1.37%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.934'
$ws.Range('E39').Value = '  +2.71%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.881'
$ws.Range('E40').Value = '  +0.96%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5071'
$ws.Range('E41').Value = '  +0.39%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1499'
$ws.Range('E42').Value = '  -1.22%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.064'
$ws.Range('E43').Value = '  +0.84%  '

$ws.Range('E44').Value = '  +0.50%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4662'
$ws.Range('E45').Value = '  +0.76%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.05'
$ws.Range('E46').Value = '  +2.27%  '

$ws.Range('E47').Value = '  -1.99%  '

$ws.Range('E48').Value = '  +0.48%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06032'

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '63.98'
$ws.Range('E50').Value = '  +0.24%  '

$ws.Range('E51').Value = '  -0.46%  '
